$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
